# Insert a new data row at row 36 (shifting the existing rows 36..118 down
# to 37..119) and populate it with the new "Poroto granado" price record for
# 2023-03-02 (Región de La Araucanía).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 36:118 down by one (xlShiftDown = -4121) to make room.
$ws.Rows("36:36").Insert(-4121)

# Populate the newly inserted row 36 with the new record.
$ws.Range("A36").Value = 10
$ws.Range("B36").Value = "Vega Modelo de Temuco"
$ws.Range("C36").Value = "La Araucanía"
$ws.Range("D36").Value = 44987
$ws.Range("E36").Value = 9
$ws.Range("F36").Value = 100112030
$ws.Range("G36").Value = "Poroto granado"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 40
$ws.Range("K36").Value = 32000
$ws.Range("L36").Value = 32000
$ws.Range("M36").Value = 32000
$ws.Range("N36").Value = "$/saco 25 kilos"
$ws.Range("O36").Value = "Región de La Araucanía"
$ws.Range("P36").Value = 1280
$ws.Range("Q36").Value = 25
$ws.Range("R36").Value = "Hortaliza"
